$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B49 text: add the "訂正" (correction) link after the original jstage link.
$ws.Range("B49").Value = "**国立医薬品食品衛生研究所**<br>[The Occurrence of Listeria monocytogenes in Imported Ready-to-Eat Foods in Japan](https://www.jstage.jst.go.jp/article/jvms/74/3/74_11-0262/_pdf/-char/en)（[訂正](https://www.jstage.jst.go.jp/article/jvms/75/2/75_11-0262e/_article/-char/ja)） <br> (Journal of Veterinary Medical Science, 2012, Volume 74, Issue 3, Pages 373-375)"

# Remove row 51 (the 2006 Miyagi Prefectural Institute of Public Health entry about
# Listeria-contaminated ready-to-eat foods) -- it was dropped from the source Markdown,
# so all subsequent rows (52-64) shift up by one.
$ws.Rows(51).Delete()
